# Add a new "Slovakia" worksheet with the Slovakia market test data.
#
# The new sheet is cloned from "Portugal" (same layout / column widths /
# merged cells / styling) and placed as the last tab, then the two
# market-specific cells are overwritten with the Slovakia values.
# Afterwards the selection / active-tab state on the touched sheets is
# updated to match what Excel leaves behind after this kind of edit.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")
$swiss = $wb.Worksheets.Item("Swiss")

# Clone Portugal into a new sheet placed right after it (use the current
# last sheet as the "After" anchor rather than $portugal itself).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Copy($null, $lastSheet)

$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Fill in the Slovakia-specific values.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3222/T3224"

# Swiss sheet's selection moves to B18.
$swiss.Range("B18").Select()

# Portugal is no longer the focused tab; its selection becomes the whole
# sheet (as if the user pressed Ctrl+A before switching off of it).
$portugal.Cells.Select()

# Slovakia (the newly added, last sheet) ends up the active / visible tab,
# selected on B2.
$slovakia.Activate()
$slovakia.Range("B2").Select()
